{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n// \"\u00a9 2020 . Contact: ...\" footer block (and the blank paragraph that\n// precedes it) that follows the \"LOB1037: \u00c0lgebra Linear (Requisito\n// fraco)\" requirement line. The blank paragraph and page-break\n// paragraph that come after the footer block are left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst marker = \"LOB1037: \u00c0lgebra Linear (Requisito fraco)\";\nconst items = paragraphs.items;\n\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex === -1) {\n  throw new Error('Could not find the \"' + marker + '\" paragraph.');\n}\n\n// The three paragraphs to drop are the ones immediately following the\n// marker: a blank paragraph, the \"Ver no Jupiter ...\" line, and the\n// \"\u00a9 2020 ...\" line.\nconst toDelete = [];\nfor (let offset = 1; offset <= 3 && markerIndex + offset < items.length; offset++) {\n  toDelete.push(items[markerIndex + offset]);\n}\n\n// Delete from the back so earlier indices/objects stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n# \"\u00a9 2020 . Contact: ...\" footer block (and the blank paragraph that\n# precedes it) that follows the \"LOB1037: \u00c0lgebra Linear (Requisito\n# fraco)\" requirement line. The blank paragraph and page-break\n# paragraph that come after the footer block are left untouched.\n\n$d = $word.ActiveDocument\n\n$marker = \"LOB1037: \u00c0lgebra Linear (Requisito fraco)\"\n$markerIndex = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (\\r),\n    # so trim before comparing.\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n    if ($text -eq $marker) {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -eq -1) {\n    throw \"Could not find the '$marker' paragraph.\"\n}\n\n# Delete the three paragraphs right after the marker (blank, \"Ver no\n# Jupiter ...\", \"\u00a9 2020 ...\") working from the back so the indices of\n# paragraphs still to be removed stay valid.\nfor ($offset = 3; $offset -ge 1; $offset--) {\n    $idx = $markerIndex + $offset\n    if ($idx -le $d.Paragraphs.Count) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
